$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSubject")

# Rename the "Maximum force [N]" column header to the new
# "Maximum contraction [Nm] / Fatigue [%]" header (with embedded line break).
# This also renames the corresponding table column automatically.
$ws.Range("E6").Value = "Maximum contraction [Nm] `n/ Fatigue [%] "

# Restore the last active cell selection as recorded after the edit.
$ws.Range("G15").Select()
